$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-27 Sunday", "2025-07-28 Monday"),
    @("774÷8=96, 6", "635÷4=158, 3"),
    @("959÷3=319, 2", "437÷4=109, 1"),
    @("700÷8=87, 4", "884÷5=176, 4"),
    @("704÷8=88, 0", "594÷7=84, 6"),
    @("200÷5=40, 0", "719÷5=143, 4"),
    @("606÷3=202, 0", "562÷3=187, 1"),
    @("120÷4=30, 0", "968÷8=121, 0"),
    @("572÷3=190, 2", "266÷3=88, 2"),
    @("265÷7=37, 6", "412÷7=58, 6"),
    @("654÷6=109, 0", "268÷4=67, 0"),
    @("292÷2=146, 0", "485÷5=97, 0"),
    @("254÷6=42, 2", "412÷3=137, 1"),
    @("368÷2=184, 0", "230÷7=32, 6"),
    @("283÷8=35, 3", "334÷9=37, 1"),
    @("703÷3=234, 1", "205÷3=68, 1"),
    @("930÷8=116, 2", "170÷5=34, 0"),
    @("408÷9=45, 3", "756÷8=94, 4"),
    @("677÷8=84, 5", "975÷7=139, 2"),
    @("425÷3=141, 2", "144÷3=48, 0"),
    @("106÷3=35, 1", "167÷6=27, 5"),
    @("521÷9=57, 8", "925÷4=231, 1"),
    @("764÷5=152, 4", "107÷6=17, 5"),
    @("174÷7=24, 6", "856÷2=428, 0"),
    @("337÷9=37, 4", "373÷3=124, 1"),
    @("392÷6=65, 2", "867÷3=289, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
